$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

# Map of cell -> new value for the "stats" sheet
$updates = @{
    "D2"  = 0.0001473743468523026
    "E2"  = 0.05991979921236634
    "G2"  = 0.003837752155959606
    "H2"  = 0.005477836821228266
    "I2"  = 0.02246203320100904
    "J2"  = 0.02479381486773491
    "K2"  = 0.0009760116226971149

    "D3"  = 0.003482725005596876
    "E3"  = 0.07466513523831964
    "G3"  = 0.003667639102786779
    "H3"  = 0.009227469563484192
    "I3"  = 0.02369260974228382
    "J3"  = 0.03439414640888572
    "K3"  = 0.0009975670836865902

    "D4"  = 0.006249985191971064
    "E4"  = 0.1103353081271052
    "G4"  = 0.00509729515761137
    "H4"  = 0.01410595281049609
    "I4"  = 0.03555978508666158
    "J4"  = 0.05052705015987158
    "K4"  = 0.001348899677395821

    "D5"  = 0.0002870620228350163
    "E5"  = 0.1013815561309457
    "G5"  = 0.005632365588098764
    "H5"  = 0.008304687682539225
    "I5"  = 0.03757016872987151
    "J5"  = 0.04413108481094241
    "K5"  = 0.001469585113227367

    "D6"  = 0.007547407876700163
    "E6"  = 0.5355360498651862
    "G6"  = 0.01110966224223375
    "H6"  = 0.02704109204933047
    "I6"  = 0.4203811264596879
    "J6"  = 0.06456113280728459
    "K6"  = 0.003652263898402452

    "D8"  = 0.0001473743468523026
    "E8"  = 0.05991979921236634
    "G8"  = 0.003837752155959606
    "H8"  = 0.005477836821228266
    "I8"  = 0.02246203320100904
    "J8"  = 0.02479381486773491
    "K8"  = 0.0009760116226971149

    "D9"  = 0.003482725005596876
    "E9"  = 0.07466513523831964
    "G9"  = 0.003667639102786779
    "H9"  = 0.009227469563484192
    "I9"  = 0.02369260974228382
    "J9"  = 0.03439414640888572
    "K9"  = 0.0009975670836865902

    "D10" = 0.006249985191971064
    "E10" = 0.1103353081271052
    "G10" = 0.00509729515761137
    "H10" = 0.01410595281049609
    "I10" = 0.03555978508666158
    "J10" = 0.05052705015987158
    "K10" = 0.001348899677395821

    "D11" = 0.0002870620228350163
    "E11" = 0.1013815561309457
    "G11" = 0.005632365588098764
    "H11" = 0.008304687682539225
    "I11" = 0.03757016872987151
    "J11" = 0.04413108481094241
    "K11" = 0.001469585113227367

    "D12" = 0.007547407876700163
    "E12" = 0.5355360498651862
    "G12" = 0.01110966224223375
    "H12" = 0.02704109204933047
    "I12" = 0.4203811264596879
    "J12" = 0.06456113280728459
    "K12" = 0.003652263898402452
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
